$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell while preserving the
# cells existing style (Excel would otherwise reformat/re-type a
# numeric-looking string and also tag the cell with a new "Text"
# number-format style when forced via a leading apostrophe).
function Set-TextValue {
    param([string]$CellRef, [string]$NewValue)
    $rng = $ws.Range($CellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $NewValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "45.132.90"
Set-TextValue "E2" "  +4.75%  "
Set-TextValue "D3" "2.432.81"
Set-TextValue "E3" "  +2.60%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.09%  "
Set-TextValue "D5" "318.16"
Set-TextValue "E5" "  +4.90%  "
Set-TextValue "D6" "104.53"
Set-TextValue "E6" "  +9.59%  "
Set-TextValue "D7" "0.518"
Set-TextValue "E7" "  +2.92%  "
Set-TextValue "E8" "  -0.09%  "
Set-TextValue "D9" "0.531"
Set-TextValue "E9" "  +10.35%  "
Set-TextValue "D10" "35.88"
Set-TextValue "E10" "  +4.41%  "
Set-TextValue "D11" "0.0805"
Set-TextValue "E11" "  +2.22%  "
Set-TextValue "E12" "  -2.50%  "
Set-TextValue "D13" "18.59"
Set-TextValue "E13" "  +2.23%  "
Set-TextValue "D14" "6.97"
Set-TextValue "E14" "  +3.02%  "
Set-TextValue "D15" "2.810.19"
Set-TextValue "E15" "  +2.58%  "
Set-TextValue "D16" "2.427.86"
Set-TextValue "E16" "  +2.75%  "
Set-TextValue "D17" "0.836"
Set-TextValue "E17" "  +4.76%  "
Set-TextValue "D18" "44.986.09"
Set-TextValue "E18" "  +4.29%  "
Set-TextValue "D19" "12.43"
Set-TextValue "E19" "  +4.06%  "
Set-TextValue "D20" "6.37"
Set-TextValue "E20" "  +1.81%  "
Set-TextValue "D21" "0.0₃0920"
Set-TextValue "E21" "  +3.71%  "
Set-TextValue "D22" "68.89"
Set-TextValue "E22" "  +1.43%  "
Set-TextValue "D23" "243.98"
Set-TextValue "E23" "  +3.70%  "
Set-TextValue "E24" "  +4.58%  "
Set-TextValue "D25" "2.50"
Set-TextValue "E25" "  +2.54%  "
Set-TextValue "E26" "  +0.06%  "
Set-TextValue "D27" "25.53"
Set-TextValue "E27" "  +4.30%  "
Set-TextValue "E28" "  -5.57%  "
Set-TextValue "D29" "9.58"
Set-TextValue "E29" "  +2.62%  "
Set-TextValue "D30" "33.78"
Set-TextValue "E30" "  +5.23%  "
Set-TextValue "D31" "48.98"
Set-TextValue "E31" "  +1.97%  "
Set-TextValue "E32" "  +17.28%  "
Set-TextValue "D33" "19.77"
Set-TextValue "E33" "  +12.46%  "
Set-TextValue "E34" "  +4.21%  "
Set-TextValue "E35" "  +0.36%  "
Set-TextValue "E36" "  +4.73%  "
Set-TextValue "D37" "1.91"
Set-TextValue "E37" "  +4.63%  "
Set-TextValue "D38" "4.53"
Set-TextValue "E38" "  +5.12%  "
Set-TextValue "D39" "127.02"
Set-TextValue "E39" "  -1.63%  "
Set-TextValue "E40" "  +1.35%  "
Set-TextValue "E41" "  +2.34%  "
Set-TextValue "E42" "  -3.17%  "
Set-TextValue "D43" "21.02"
Set-TextValue "E43" "  +0.16%  "
Set-TextValue "E44" "  +4.57%  "
Set-TextValue "D45" "1.944.83"
Set-TextValue "E45" "  +0.74%  "
Set-TextValue "E46" "  -0.55%  "
Set-TextValue "D47" "2.96"
Set-TextValue "E47" "  +8.64%  "
Set-TextValue "D48" "9.22"
Set-TextValue "E48" "  -0.27%  "
Set-TextValue "D49" "1.79"
Set-TextValue "E49" "  +18.90%  "
Set-TextValue "D50" "75.99"
Set-TextValue "E50" "  +6.46%  "
Set-TextValue "D51" "54.32"
Set-TextValue "E51" "  +5.60%  "
